$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MSG: None`n`nMSG: I have recorded the decision as no_decision.`n"
$ws.Range("C3").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been successfully made.`n"
$ws.Range("C4").Value = "MSG: None`n`nMSG: The decision on which movie to show on Friday was not reached, and therefore no action can be taken at this time.`n"
$ws.Range("C6").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights to both movies.`n"
$ws.Range("C7").Value = "MSG: None`n`nMSG: The decision has been recorded: no movie will be shown on Friday as the committee did not reach a conclusion.`n"
$ws.Range("D7").Value = "no_decision, "
$ws.Range("C8").Value = "MSG: None`n`nMSG: The rights to both movies have been acquired successfully.`n"
$ws.Range("D8").Value = "both_movies, "
$ws.Range("C9").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been successfully recorded.`n"
$ws.Range("C10").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie`" to be shown on Friday.`n"
$ws.Range("C11").Value = "MSG: None`n`nMSG: The decision has been recorded indicating that no agreement was reached on which movie to show on Friday.`n"
$ws.Range("C12").Value = "MSG: None`n`nMSG: The decision has been recorded. No agreement was reached regarding the movie to be shown on Friday.`n"
$ws.Range("C13").Value = "MSG: None`n`nMSG: The decision to show `"Barbie`" has been successfully recorded.`n"
$ws.Range("C14").Value = "MSG: None`n`nMSG: The decision has been recorded: no movie was selected for Friday.`n"
$ws.Range("C15").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Oppenheimer`" as the movie to be shown on Friday.`n"
$ws.Range("C16").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights to `"Barbie`" for the showing on Friday.`n"
$ws.Range("D16").Value = "Barbie_was_selected, "
$ws.Range("C17").Value = "MSG: None`n`nMSG: The decision regarding the movie to show on Friday has ended without a specific choice being made.`n"
$ws.Range("C18").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been made.`n"
$ws.Range("C19").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday did not reach a consensus, resulting in no decision being made.`n"
$ws.Range("C20").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights to `"Barbie.`"`n"
$ws.Range("C21").Value = "MSG: None`n`nMSG: The decision to acquire the rights for both movies has been successfully recorded.`n"
$ws.Range("C22").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie being selected for Friday.`n"
$ws.Range("C23").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no movie was chosen for Friday.`n"
$ws.Range("C24").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no movie will be shown on Friday.`n"
$ws.Range("C25").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" to be shown on Friday.`n"
$ws.Range("C26").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie`" to be shown on Friday.`n"
$ws.Range("C27").Value = "MSG: None`n`nMSG: No decision was made regarding the movie to show on Friday.`n"
$ws.Range("C28").Value = "MSG: None`n`nMSG: The decision regarding the movie to be shown on Friday resulted in no agreement, hence no decision has been made.`n"
$ws.Range("C29").Value = "MSG: None`n`nMSG: The decision has been recorded successfully. The movie `"Barbie`" will be shown on Friday.`n"
$ws.Range("C30").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" as the movie to be shown on Friday.`n"
$ws.Range("C31").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no_decision.`" The committee did not reach a conclusion on which movie to show on Friday.`n"
$ws.Range("C32").Value = "MSG: None`n`nMSG: The decision has been recorded with no agreement on which movie to show on Friday.`n"
$ws.Range("D32").Value = "no_decision, "
$ws.Range("C33").Value = "MSG: None`n`nMSG: The decision has been recorded with no choice made regarding the movie for Friday.`n"
$ws.Range("C34").Value = "MSG: None`n`nMSG: The decision has been recorded, and the movie `"Barbie`" will be acquired for Friday's showing.`n"
$ws.Range("C35").Value = "MSG: None`n`nMSG: The conversation did not lead to a decision about which movie to show on Friday.`n"
$ws.Range("C36").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision was made regarding Friday’s movie.`n"
$ws.Range("C37").Value = "MSG: None`n`nMSG: The decision has been recorded, and `"Barbie`" has been selected to be shown on Friday.`n"
$ws.Range("C38").Value = "MSG: None`n`nMSG: The decision process concluded without a selection for Friday's movie, indicating that no agreement was reached.`n"
$ws.Range("C39").Value = "MSG: None`n`nMSG: I have successfully recorded the decision to acquire the rights for both movies.`n"
$ws.Range("C40").Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday was not finalized, so no action has been taken.`n"
$ws.Range("C41").Value = "MSG: None`n`nMSG: The decision has been recorded, and rights for `"Barbie`" will be acquired.`n"
$ws.Range("C42").Value = "MSG: None`n`nMSG: The decision has been recorded, and there is no definitive plan for what movie will be shown on Friday.`n"
$ws.Range("D42").Value = "no_decision, "
$ws.Range("C43").Value = "MSG: None`n`nMSG: The decision has been recorded with the function indicating no consensus was reached regarding the movie for Friday.`n"
$ws.Range("C44").Value = "MSG: None`n`nMSG: The committee reached no decision regarding the movie to be shown on Friday.`n"
$ws.Range("C45").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie being selected for Friday.`n"
$ws.Range("C46").Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday has not been made.`n"
$ws.Range("C47").Value = "MSG: None`n`nMSG: The decision regarding the movie for Friday has not been made, and thus there is no selection to proceed with.`n"
$ws.Range("C48").Value = "MSG: None`n`nMSG: No decision was made about the movie to show on Friday.`n"
$ws.Range("C49").Value = "MSG: None`n`nMSG: The decision was made to conclude the conversation without selecting a movie for Friday.`n"
$ws.Range("C50").Value = "MSG: None`n`nMSG: The decision process concluded without selecting a movie, and thus the no_deciison function has been called successfully.`n"
$ws.Range("C51").Value = "MSG: None`n`nMSG: The rights to both movies have been acquired successfully.`n"
$ws.Range("C52").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie was selected for Friday.`n"
$ws.Range("C53").Value = "MSG: None`n`nMSG: I have recorded the decision to acquire the rights for `"Barbie`" as the movie to be shown on Friday.`n"
$ws.Range("C54").Value = "MSG: None`n`nMSG: The rights for both `"Barbie`" and `"Oppenheimer`" have been successfully acquired.`n"
$ws.Range("D54").Value = "both_movies, "
$ws.Range("C55").Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday was not finalized, resulting in no selection being made.`n"
$ws.Range("D55").Value = "no_decision, "
$ws.Range("C56").Value = "MSG: None`n`nMSG: The decision on which movie to show on Friday has not been made.`n"
$ws.Range("C57").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Oppenheimer`" has been successfully recorded.`n"
$ws.Range("C58").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been recorded.`n"
$ws.Range("C59").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision`" regarding the movie selection for Friday.`n"
$ws.Range("C60").Value = "MSG: None`n`nMSG: The decision process concluded without selecting a movie, so no acquisition will be made.`n"
$ws.Range("C61").Value = "MSG: None`n`nMSG: The decision has been recorded successfully to acquire the rights for `"Barbie`" to be shown on Friday.`n"
$ws.Range("C62").Value = "MSG: None`n`nMSG: The decision to acquire the rights for both movies has been made.`n"
$ws.Range("C63").Value = "MSG: None`n`nMSG: The decision has been recorded as no choice of a movie for Friday.`n"
$ws.Range("C64").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie was selected.`n"
$ws.Range("C65").Value = "MSG: None`n`nMSG: The decision process concluded without reaching a consensus on which movie to show on Friday. Therefore, no decision was made.`n"
$ws.Range("C66").Value = "MSG: None`n`nMSG: The decision about what movie to show on Friday has not been made.`n"
$ws.Range("C67").Value = "MSG: None`n`nMSG: The decision about the movie to show on Friday remains unresolved.`n"
$ws.Range("C68").Value = "MSG: None`n`nMSG: I have recorded the decision as `"no_decision`" based on the guidelines provided and the lack of a clear agreement on the movie for Friday.`n"
